{"js": "// Apply the four text edits described in the commit:\n//  1. \"CHPC instructions\" -> \"CHPC instructions, done by Mike\"\n//  2. \"Implemented by Mike.\" (Parallel Shared Memory CPU section)\n//       -> \"Implemented by Mike. Straightforward to implement off of the\n//           back of the serial implementation.\"\n//  3. \"Implemented by Mike.\" (Parallel CUDA GPU section)\n//       -> \"Implemented by Mike. Some refactoring was necessary here to\n//           transfer all data structures away from vectors (which are not\n//           supported by kernel functions).\"\n//  4. \"...processing across MPI ranks.\" (Distributed Memory GPU section)\n//       -> \"...processing across MPI ranks. He also handled all of the\n//           CHPC integration.\"\n\nconst body = context.document.body;\n\n// --- Edit 1: unique text, simple search + replace ------------------------\nconst chpcResults = body.search(\"CHPC instructions\", { matchCase: true });\nchpcResults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < chpcResults.items.length; i++) {\n    const r = chpcResults.items[i];\n    if (r.text === \"CHPC instructions\") {\n        r.insertText(\"CHPC instructions, done by Mike\", Word.InsertLocation.replace);\n    }\n}\nawait context.sync();\n\n// --- Edits 2 & 3: \"Implemented by Mike.\" appears twice (identical text) --\n// Search returns them in document order, so item 0 is the \"Parallel Shared\n// Memory CPU\" paragraph and item 1 is the \"Parallel CUDA GPU\" paragraph.\nconst implResults = body.search(\"Implemented by Mike.\", { matchCase: true });\nimplResults.load(\"items,text\");\nawait context.sync();\n\nconst implReplacements = [\n    \"Implemented by Mike. Straightforward to implement off of the back of the serial implementation.\",\n    \"Implemented by Mike. Some refactoring was necessary here to transfer all data structures away from vectors (which are not supported by kernel functions).\"\n];\n\nlet implIdx = 0;\nfor (let i = 0; i < implResults.items.length; i++) {\n    const r = implResults.items[i];\n    if (r.text === \"Implemented by Mike.\" && implIdx < implReplacements.length) {\n        r.insertText(implReplacements[implIdx], Word.InsertLocation.replace);\n        implIdx++;\n    }\n}\nawait context.sync();\n\n// --- Edit 4: append sentence onto the end of the MPI-ranks paragraph -----\n// Using InsertLocation.replace with the full target substring (rather than\n// InsertLocation.after) keeps the edit inside the existing run so its\n// formatting/rsid attributes are preserved, matching how Office.js handled\n// edits 1-3 above.\nconst mpiResults = body.search(\"processing across MPI ranks.\", { matchCase: true });\nmpiResults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < mpiResults.items.length; i++) {\n    const r = mpiResults.items[i];\n    if (r.text === \"processing across MPI ranks.\") {\n        r.insertText(\n            \"processing across MPI ranks. He also handled all of the CHPC integration.\",\n            Word.InsertLocation.replace\n        );\n    }\n}\nawait context.sync();\n", "ps1": "# Apply the four text edits described in the commit:\n#  1. \"CHPC instructions\" -> \"CHPC instructions, done by Mike\"\n#  2. \"Implemented by Mike.\" (Parallel Shared Memory CPU section)\n#       -> \"Implemented by Mike. Straightforward to implement off of the\n#           back of the serial implementation.\"\n#  3. \"Implemented by Mike.\" (Parallel CUDA GPU section)\n#       -> \"Implemented by Mike. Some refactoring was necessary here to\n#           transfer all data structures away from vectors (which are not\n#           supported by kernel functions).\"\n#  4. \"...processing across MPI ranks.\" (Distributed Memory GPU section)\n#       -> \"...processing across MPI ranks. He also handled all of the\n#           CHPC integration.\"\n#\n# InsertAfter() on the Find-matched range appends text onto the end of the\n# existing run (keeping its formatting / rsid attributes intact) rather than\n# minting a brand-new run, which mirrors how the target document's XML\n# changed (same <w:r>, extended <w:t>).\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: unique text, find + append ----------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"CHPC instructions\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nif ($rng.Find.Execute()) {\n    $rng.InsertAfter(\", done by Mike\")\n}\n\n# --- Edits 2 & 3: \"Implemented by Mike.\" appears twice (identical text) --\n# Find returns matches in document order, so the first hit is the \"Parallel\n# Shared Memory CPU\" paragraph and the second is \"Parallel CUDA GPU\".\n$implAppends = @(\n    \" Straightforward to implement off of the back of the serial implementation.\",\n    \" Some refactoring was necessary here to transfer all data structures away from vectors (which are not supported by kernel functions).\"\n)\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Implemented by Mike.\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\n$i = 0\nwhile ($i -lt $implAppends.Count -and $rng.Find.Execute()) {\n    $rng.InsertAfter($implAppends[$i])\n    $i++\n    $rng.Collapse(0)\n}\n\n# --- Edit 4: append sentence onto the end of the MPI-ranks paragraph -----\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"processing across MPI ranks.\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nif ($rng.Find.Execute()) {\n    $rng.InsertAfter(\" He also handled all of the CHPC integration.\")\n}\n"}
